$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint2")
$ws.Range("A24").Value = "Test"
